$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''305.34'
$ws.Range("E2").Value = '''-0.53%'
$ws.Range("D3").Value = '''35.77'
$ws.Range("E3").Value = '''-0.39%'
$ws.Range("D4").Value = '''5.038'
$ws.Range("E4").Value = '''-1.49%'
$ws.Range("D5").Value = '''0.07964'
$ws.Range("E5").Value = '''-1.42%'
$ws.Range("D6").Value = '''1.910'
$ws.Range("E6").Value = '''-1.89%'
$ws.Range("B7").Value = 'KuCoinToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D7").Value = '''7.770'
$ws.Range("E7").Value = '''0.24%'
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = '''0.9201'
$ws.Range("E8").Value = '''-0.83%'
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").Value = '''0.1299'
$ws.Range("E9").Value = '''-5.27%'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '''0.1911'
$ws.Range("E10").Value = '''0.10%'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = '''0.09114'
$ws.Range("E11").Value = '''-0.64%'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = '''0.03417'
$ws.Range("E12").Value = '''-1.10%'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = '''0.09843'
$ws.Range("E13").Value = '''0.14%'
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = '''0.001407'
$ws.Range("E14").Value = '''-0.26%'
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = '''0.006182'
$ws.Range("E15").Value = '''6.95%'
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").Value = '''3.718'
$ws.Range("E16").Value = '''2.74%'
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").Value = '''4.135'
$ws.Range("E17").Value = '''-1.52%'
$ws.Range("D18").Value = '''3.395'
$ws.Range("E18").Value = '''12.94%'
$ws.Range("D19").Value = '''0.3445'
$ws.Range("E19").Value = '''-0.05%'
$ws.Range("D20").Value = '''0.1311'
$ws.Range("E20").Value = '''-2.25%'
$ws.Range("D21").Value = '''5.161'
$ws.Range("E21").Value = '''5.13%'
$ws.Range("D23").Value = '''0.04408'
$ws.Range("E23").Value = '''-0.71%'
$ws.Range("D24").Value = '''0.001233'
$ws.Range("E24").Value = '''1.04%'
$ws.Range("D25").Value = '''0.004632'
$ws.Range("E25").Value = '''-4.17%'
$ws.Range("D26").Value = '''0.0001250'
$ws.Range("E26").Value = '''0.47%'
$ws.Range("D27").Value = '''0.0004442'
$ws.Range("E27").Value = '''0.01%'
$ws.Range("D39").Value = '''0.01942'
$ws.Range("E39").Value = '''-4.03%'
$ws.Range("D40").Value = '''0.05275'
$ws.Range("E40").Value = '''7.18%'
$ws.Range("D41").Value = '''0.007588'
$ws.Range("E41").Value = '''-1.56%'
$ws.Range("D42").Value = '''0.01012'
$ws.Range("E42").Value = '''0.03%'
$ws.Range("D43").Value = '''0.1351'
$ws.Range("E43").Value = '''-1.96%'
$ws.Range("D44").Value = '''0.002159'
$ws.Range("E44").Value = '''2.51%'
$ws.Range("D45").Value = '''0.009941'
$ws.Range("E45").Value = '''-14.28%'
$ws.Range("D46").Value = '''0.00006098'
$ws.Range("E46").Value = '''-5.71%'
$ws.Range("D47").Value = '''0.00000000750'
$ws.Range("E47").Value = '''-0.28%'
$ws.Range("D48").Value = '''65.22'
$ws.Range("E48").Value = '''2.60%'
$ws.Range("D49").Value = '''0.001658'
$ws.Range("E49").Value = '''39.01%'
$ws.Range("D50").Value = '''0.00002101'
$ws.Range("E50").Value = '''-0.28%'
$ws.Range("D51").Value = '''0.0002001'
$ws.Range("E51").Value = '''-0.28%'
